$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 45220
$ws.Range("D2").Value = 45219
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 45227
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 45226
$ws.Range("C5").Value = 8
$ws.Range("D5").Value = 45227
$ws.Range("D6").Value = 45219
$ws.Range("D7").Value = 45219
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 45221
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 45235
$ws.Range("C10").Value = 8
$ws.Range("D10").Value = 45226
$ws.Range("C11").Value = 8
$ws.Range("D11").Value = 45227
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = 45221
$ws.Range("C13").Value = 8
$ws.Range("D13").Value = 45226
$ws.Range("C15").Value = 8
$ws.Range("D15").Value = 45225
$ws.Range("D17").Value = 45220
$ws.Range("D19").Value = 45219
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 45227
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 45233
$ws.Range("D22").Value = 45219
$ws.Range("C23").Value = 8
$ws.Range("D23").Value = 45226
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 45227
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 45226
$ws.Range("D27").Value = 45219
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = 45223
$ws.Range("D30").Value = 45219
$ws.Range("C31").Value = 16
$ws.Range("D31").Value = 45233
$ws.Range("C32").Value = 16
$ws.Range("D32").Value = 45233
$ws.Range("D33").Value = 45219
$ws.Range("C34").Value = 8
$ws.Range("D34").Value = 45225
$ws.Range("C35").Value = 8
$ws.Range("D35").Value = 45227
$ws.Range("D36").Value = 45219
$ws.Range("C38").Value = 4
$ws.Range("D38").Value = 45222
$ws.Range("C39").Value = 2
$ws.Range("D39").Value = 45221
$ws.Range("D40").Value = 45219
$ws.Range("C41").Value = 8
$ws.Range("D41").Value = 45225
$ws.Range("C42").Value = 4
$ws.Range("D42").Value = 45222
$ws.Range("C43").Value = 2
$ws.Range("D43").Value = 45220
$ws.Range("C44").Value = 16
$ws.Range("D44").Value = 45234
$ws.Range("C45").Value = 8
$ws.Range("D45").Value = 45226
$ws.Range("C46").Value = 16
$ws.Range("D46").Value = 45235
$ws.Range("C48").Value = 4
$ws.Range("D48").Value = 45222
$ws.Range("D49").Value = 45219
$ws.Range("C50").Value = 2
$ws.Range("D50").Value = 45220
$ws.Range("C51").Value = 2
$ws.Range("D51").Value = 45220
$ws.Range("C52").Value = 4
$ws.Range("D52").Value = 45222
$ws.Range("C54").Value = 8
$ws.Range("D54").Value = 45225
$ws.Range("D55").Value = 45220
$ws.Range("C56").Value = 4
$ws.Range("D56").Value = 45222
$ws.Range("C57").Value = 16
$ws.Range("D57").Value = 45234
$ws.Range("C58").Value = 16
$ws.Range("D58").Value = 45234
$ws.Range("D59").Value = 45219
$ws.Range("D61").Value = 45219
$ws.Range("C63").Value = 8
$ws.Range("D63").Value = 45226
$ws.Range("C64").Value = 4
$ws.Range("D64").Value = 45221
$ws.Range("D65").Value = 45219
$ws.Range("C66").Value = 4
$ws.Range("D66").Value = 45222
$ws.Range("C68").Value = 4
$ws.Range("D68").Value = 45221
$ws.Range("C69").Value = 4
$ws.Range("D69").Value = 45222
$ws.Range("D70").Value = 45219
$ws.Range("C71").Value = 4
$ws.Range("D71").Value = 45221
$ws.Range("D73").Value = 45219
$ws.Range("C74").Value = 2
$ws.Range("D74").Value = 45220
$ws.Range("C75").Value = 2
$ws.Range("D75").Value = 45221
$ws.Range("D76").Value = 45219
$ws.Range("C77").Value = 8
$ws.Range("D77").Value = 45227
$ws.Range("C78").Value = 4
$ws.Range("D78").Value = 45223
$ws.Range("C79").Value = 8
$ws.Range("D79").Value = 45226
$ws.Range("D80").Value = 45219
$ws.Range("C81").Value = 8
$ws.Range("D81").Value = 45226
$ws.Range("C82").Value = 8
$ws.Range("D82").Value = 45226
$ws.Range("D83").Value = 45219
$ws.Range("C84").Value = 2
$ws.Range("D84").Value = 45221
$ws.Range("D85").Value = 45219
$ws.Range("D86").Value = 45219
$ws.Range("C89").Value = 2
$ws.Range("D89").Value = 45220
$ws.Range("C90").Value = 8
$ws.Range("D90").Value = 45227
$ws.Range("C91").Value = 8
$ws.Range("D91").Value = 45226
$ws.Range("C92").Value = 2
$ws.Range("D92").Value = 45221
$ws.Range("C93").Value = 8
$ws.Range("D93").Value = 45227
$ws.Range("C94").Value = 4
$ws.Range("D94").Value = 45221
$ws.Range("C95").Value = 2
$ws.Range("D95").Value = 45220
$ws.Range("D96").Value = 45219
$ws.Range("C97").Value = 2
$ws.Range("D97").Value = 45220
$ws.Range("C98").Value = 4
$ws.Range("D98").Value = 45223
$ws.Range("D99").Value = 45219
$ws.Range("C100").Value = 4
$ws.Range("D100").Value = 45223
$ws.Range("D101").Value = 45219
$ws.Range("C102").Value = 2
$ws.Range("D102").Value = 45221
$ws.Range("D103").Value = 45220
$ws.Range("C104").Value = 4
$ws.Range("D104").Value = 45221
$ws.Range("D105").Value = 45219
$ws.Range("C106").Value = 2
$ws.Range("D106").Value = 45220
$ws.Range("C107").Value = 16
$ws.Range("D107").Value = 45234
$ws.Range("D108").Value = 45219
$ws.Range("D109").Value = 45219
$ws.Range("C110").Value = 2
$ws.Range("D110").Value = 45220
$ws.Range("D111").Value = 45219
$ws.Range("C112").Value = 8
$ws.Range("D112").Value = 45226
$ws.Range("D113").Value = 45219
$ws.Range("C115").Value = 2
$ws.Range("D115").Value = 45221
$ws.Range("C116").Value = 2
$ws.Range("D116").Value = 45220
$ws.Range("C117").Value = 2
$ws.Range("D117").Value = 45220
$ws.Range("C118").Value = 4
$ws.Range("D118").Value = 45222
$ws.Range("D119").Value = 45220
$ws.Range("C120").Value = 16
$ws.Range("D120").Value = 45235
$ws.Range("D121").Value = 45219
$ws.Range("C122").Value = 8
$ws.Range("D122").Value = 45226
$ws.Range("C123").Value = 4
$ws.Range("D123").Value = 45221
$ws.Range("D124").Value = 45219
$ws.Range("C125").Value = 2
$ws.Range("D125").Value = 45220
$ws.Range("C126").Value = 16
$ws.Range("D126").Value = 45234
$ws.Range("C127").Value = 8
$ws.Range("D127").Value = 45227
$ws.Range("C128").Value = 4
$ws.Range("D128").Value = 45223
$ws.Range("C129").Value = 2
$ws.Range("D129").Value = 45220
$ws.Range("C130").Value = 8
$ws.Range("D130").Value = 45226
$ws.Range("C131").Value = 2
$ws.Range("D131").Value = 45220
